# "Updated PID Code for 0610"
#
# 1. Update the four measured input values that feed the PID-tuning
#    formulas (columns B/D, rows 34-35). All the downstream formulas in
#    rows 39-49 recompute automatically from these.
# 2. Remove the scratch "Time"/"PID" scatter-chart data block (A52:B66)
#    and its chart object - it was a temporary plot no longer needed.
# 3. Remove the now-unused "_xlchart.*" defined names that backed that
#    chart's series references.
# 4. Restore the sheet view (zoom + selection) to match the saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Updated measurements ---------------------------------------------
$ws.Range("B34").Value = 90.51
$ws.Range("D34").Value = 73
$ws.Range("B35").Value = 109
$ws.Range("D35").Value = 105

# --- 2. Drop the scratch "Time" / "PID" chart + its backing data ---------
$ws.ChartObjects("Chart 5").Delete()

$ws.Range("A52:B52").ClearContents()
$ws.Range("A53:D66").ClearContents()

# --- 3. Remove the leftover _xlchart defined names ------------------------
$nameCount = $wb.Names.Count
if ($nameCount -gt 0) {
    $namesToDelete = @()
    for ($i = 1; $i -le $nameCount; $i++) {
        $namesToDelete += $wb.Names.Item($i).Name
    }
    foreach ($n in $namesToDelete) {
        $wb.Names.Item($n).Delete()
    }
}

# --- 4. Sheet view: zoom to 92% and select B40 -----------------------------
$excel.ActiveWindow.Zoom = 92
$ws.Range("B40").Select()
